$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting (values look numeric but must stay as text)
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "43.014.56"
$ws.Range("E2").Value = "  +2.97%  "

# Row 3
$ws.Range("D3").Value = "2.293.73"
$ws.Range("E3").Value = "  +1.77%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "310.17"
$ws.Range("E5").Value = "  +2.03%  "

# Row 6
$ws.Range("D6").Value = "100.74"
$ws.Range("E6").Value = "  +7.13%  "

# Row 7
$ws.Range("E7").Value = "  +2.57%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("E9").Value = "  +7.54%  "

# Row 10
$ws.Range("D10").Value = "35.90"
$ws.Range("E10").Value = "  +3.78%  "

# Row 11
$ws.Range("D11").Value = "0.0821"
$ws.Range("E11").Value = "  +4.48%  "

# Row 12
$ws.Range("E12").Value = "  +0.88%  "

# Row 13
$ws.Range("D13").Value = "7.10"
$ws.Range("E13").Value = "  +7.64%  "

# Row 14
$ws.Range("D14").Value = "2.652.66"
$ws.Range("E14").Value = "  +1.92%  "

# Row 15
$ws.Range("D15").Value = "14.96"
$ws.Range("E15").Value = "  +4.75%  "

# Row 16
$ws.Range("D16").Value = "2.300.98"
$ws.Range("E16").Value = "  -1.03%  "

# Row 17
$ws.Range("D17").Value = "0.806"
$ws.Range("E17").Value = "  +2.68%  "

# Row 18
$ws.Range("D18").Value = "42.962.94"
$ws.Range("E18").Value = "  +3.14%  "

# Row 19
$ws.Range("D19").Value = "12.47"
$ws.Range("E19").Value = "  +1.72%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0924"
$ws.Range("E20").Value = "  +3.17%  "

# Row 21
$ws.Range("D21").Value = "6.07"
$ws.Range("E21").Value = "  +2.29%  "

# Row 22
$ws.Range("D22").Value = "68.40"
$ws.Range("E22").Value = "  +0.91%  "

# Row 23
$ws.Range("D23").Value = "240.03"
$ws.Range("E23").Value = "  +1.71%  "

# Row 24
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "2.62"
$ws.Range("E24").Value = "  +2.89%  "

# Row 25
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "2.00"
$ws.Range("E25").Value = "  +4.79%  "

# Row 26
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("D27").Value = "24.41"
$ws.Range("E27").Value = "  +3.67%  "

# Row 28
$ws.Range("D28").Value = "38.37"
$ws.Range("E28").Value = "  +6.69%  "

# Row 29
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  +9.53%  "

# Row 30
$ws.Range("D30").Value = "9.63"
$ws.Range("E30").Value = "  +2.40%  "

# Row 31
$ws.Range("D31").Value = "167.07"
$ws.Range("E31").Value = "  +4.83%  "

# Row 32
$ws.Range("D32").Value = "5.30"
$ws.Range("E32").Value = "  +2.30%  "

# Row 33
$ws.Range("E33").Value = "  +0.14%  "

# Row 34
$ws.Range("E34").Value = "  -0.59%  "

# Row 35
$ws.Range("D35").Value = "17.68"
$ws.Range("E35").Value = "  +4.71%  "

# Row 36
$ws.Range("D36").Value = "0.0739"
$ws.Range("E36").Value = "  +0.97%  "

# Row 37
$ws.Range("E37").Value = "  +3.04%  "

# Row 39
$ws.Range("E39").Value = "  +2.18%  "

# Row 40
$ws.Range("E40").Value = "  +1.09%  "

# Row 41
$ws.Range("D41").Value = "4.21"
$ws.Range("E41").Value = "  +6.49%  "

# Row 42
$ws.Range("D42").Value = "2.29"
$ws.Range("E42").Value = "  -0.28%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0288"
$ws.Range("E43").Value = "  +2.90%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.969.40"
$ws.Range("E44").Value = "  +0.50%  "

# Row 45
$ws.Range("D45").Value = "19.10"
$ws.Range("E45").Value = "  +2.48%  "

# Row 46
$ws.Range("D46").Value = "3.02"
$ws.Range("E46").Value = "  +3.91%  "

# Row 47
$ws.Range("D47").Value = "9.84"
$ws.Range("E47").Value = "  +0.24%  "

# Row 48
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "2.95"
$ws.Range("E48").Value = "  +17.95%  "

# Row 49
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "55.64"
$ws.Range("E49").Value = "  +5.78%  "

# Row 50
$ws.Range("D50").Value = "2.522.44"
$ws.Range("E50").Value = "  +1.89%  "

# Row 51
$ws.Range("E51").Value = "  +2.42%  "

# Restore default style on column D (remove temporary text-format style index)
$colD.Style = "Normal"
